$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = 51
$ws.Cells.Item(40, 3).Value = 9
$ws.Cells.Item(40, 4).Value = 14
$ws.Cells.Item(40, 5).Value = 17
$ws.Cells.Item(40, 6).Value = 74
$ws.Cells.Item(40, 7).Value = 91

$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = 51
$ws.Cells.Item(41, 3).Value = 9
$ws.Cells.Item(41, 4).Value = 14
$ws.Cells.Item(41, 5).Value = 17
$ws.Cells.Item(41, 6).Value = 74
$ws.Cells.Item(41, 7).Value = 91
